$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: drop the two "category" columns, keep Name / Email ID / PhoneNo ---

# Fix the stray space in the e-mail address before the columns shift around
$ws1.Range("D2").Value = "abc@gmail.com"

# Remove columns A:B (category / Category) - C:F shift left into A:D
$ws1.Columns("A:B").Delete()

# The custom font that forced this row's height is gone now; drop back to default
$ws1.Rows("3").AutoFit()

# The hyperlinks don't automatically follow the column shift, so recreate them
# at their new home cells (B2, B3)
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:abc@gmail.com")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:John@gmail.com")

# Carry the formatted area one row further down
$ws1.Range("B6").Style = "Hyperlink"

$ws1.Range("C6").Select()

# --- Add Sheet2: Email_ID / PhoneNumber summary ---

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Email_ID"
$ws2.Range("B1").Value = "PhoneNumber"
$ws2.Range("A2").Value = "abc@gmail.com"
$ws2.Range("B2").Value = 9390799743
$ws2.Range("A3").Value = "John@gmail.com"
$ws2.Range("B3").Value = 1235467890

$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:abc@gmail.com")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:John@gmail.com")

$ws2.Columns("A").ColumnWidth = 16.7265625
$ws2.Columns("B").ColumnWidth = 10.81640625

$ws2.Activate()
